$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-06 Monday" "2025-10-07 Tuesday"

Replace-Text "709×6=" "283×5="
Replace-Text "556×6=" "669×2="
Replace-Text "733×5=" "433×8="
Replace-Text "258×3=" "617×4="
Replace-Text "464×2=" "576×9="
Replace-Text "559×7=" "636×3="
Replace-Text "234×3=" "787×5="
Replace-Text "706×7=" "898×4="
Replace-Text "749×2=" "807×4="
Replace-Text "695×5=" "511×8="
Replace-Text "668×4=" "232×2="
Replace-Text "745×3=" "689×8="
Replace-Text "347×2=" "562×5="
Replace-Text "444×7=" "523×3="
Replace-Text "474×8=" "292×7="
Replace-Text "139×6=" "495×4="
Replace-Text "922×8=" "698×3="
Replace-Text "816×7=" "774×6="
Replace-Text "139×9=" "128×4="
Replace-Text "359×4=" "950×9="
Replace-Text "593×9=" "544×6="
Replace-Text "572×9=" "293×3="
Replace-Text "597×9=" "389×8="
Replace-Text "712×2=" "298×9="
Replace-Text "128×8=" "656×3="
